$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 1

$ws.Range("B3").Value = 0
$ws.Range("D3").Value = 5
$ws.Range("E3").Value = 2

$ws.Range("C4").Value = 7
$ws.Range("E4").Value = 16

$ws.Range("C5").Value = 3
$ws.Range("D5").Value = 21
$ws.Range("F5").Value = 10

$ws.Range("E6").Value = 2
$ws.Range("G6").Value = 0

$ws.Range("F7").Value = 2
